$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "price" column (C) values for each disease row, 1:1 with rows 1-19.
$prices = @("15.0","25.0","12.0","34.0","23.0","12.0","11.0","23.0","24.0","27.0","12.0","18.0","19.0","29.0","38.0","34.0","25.0","38.0","37.0")

$target = $ws.Range("C1:C19")

# Temporarily force the column to Text so the numeric-looking price strings
# ("15.0", "25.0", ...) are stored as literal text rather than being
# auto-converted to numbers when assigned, then clear the formatting again so
# the cells end up with the workbook's default (unformatted) style.
$target.NumberFormat = "@"

for ($i = 1; $i -le 19; $i++) {
    $ws.Range("C$i").Value = $prices[$i - 1]
}

$target.ClearFormats()

$ws.Range("C20").Select() | Out-Null
